# Auto-generated edit script: update per-row H/I/J/K/L/M/N financial values
# across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1499.75
$ws.Range("J32").Value = 2399.3333
$ws.Range("L32").Value = 2399.3333
$ws.Range("N32").Value = -3051.3333
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = $null
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = $null
$ws.Range("H76").Value = 5262.2104
$ws.Range("I76").Value = 3657.1667
$ws.Range("K76").Value = 3657.1667
$ws.Range("M76").Value = -3342.1667
$ws.Range("H79").Value = 5262.2104
$ws.Range("I79").Value = 3657.1667
$ws.Range("K79").Value = 3657.1667
$ws.Range("M79").Value = -2565.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4101.8887
$ws.Range("I32").Value = 3967.8408
$ws.Range("K32").Value = 3967.8408
$ws.Range("M32").Value = -3680.8408
$ws.Range("H63").Value = 1388.4
$ws.Range("I63").Value = 1388.4
$ws.Range("K63").Value = 1388.4
$ws.Range("M63").Value = -702.4000000000001
$ws.Range("H66").Value = 1388.4
$ws.Range("I66").Value = 1388.4
$ws.Range("K66").Value = 6942
$ws.Range("M66").Value = -3510
$ws.Range("H74").Value = 4881
$ws.Range("I74").Value = 1731.4615
$ws.Range("K74").Value = 1731.4615
$ws.Range("M74").Value = -857.4614999999999
$ws.Range("H77").Value = 4881
$ws.Range("I77").Value = 1731.4615
$ws.Range("K77").Value = 8657.307499999999
$ws.Range("M77").Value = -4289.307499999999
$ws.Range("H122").Value = 314920.3
$ws.Range("I122").Value = 2820.7646
$ws.Range("K122").Value = 8462.293799999999
$ws.Range("M122").Value = -6012.293799999999
$ws.Range("H132").Value = 7055.143
$ws.Range("I132").Value = 7812.1304
$ws.Range("K132").Value = 23436.3912
$ws.Range("M132").Value = -20906.3912

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7790.2666
$ws.Range("I86").Value = 12056.75
$ws.Range("J86").Value = 2914.2856
$ws.Range("K86").Value = 12056.75
$ws.Range("L86").Value = 2914.2856
$ws.Range("M86").Value = -10933.75
$ws.Range("N86").Value = -5160.2856
$ws.Range("H89").Value = 7790.2666
$ws.Range("I89").Value = 12056.75
$ws.Range("J89").Value = 2914.2856
$ws.Range("K89").Value = 60283.75
$ws.Range("L89").Value = 14571.428
$ws.Range("M89").Value = -54667.75
$ws.Range("N89").Value = -25803.428
$ws.Range("H105").Value = 95727
$ws.Range("I105").Value = 127374.625
$ws.Range("K105").Value = 127374.625
$ws.Range("M105").Value = -125627.625
$ws.Range("H134").Value = 10490.275
$ws.Range("J134").Value = 3375
$ws.Range("L134").Value = 10125
$ws.Range("N134").Value = -15195

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16500.555
$ws.Range("J31").Value = 3668
$ws.Range("L31").Value = 3668
$ws.Range("N31").Value = -4258
$ws.Range("H34").Value = 16500.555
$ws.Range("J34").Value = 3668
$ws.Range("L34").Value = 3668
$ws.Range("N34").Value = -4072
$ws.Range("H94").Value = 2930.4546
$ws.Range("J94").Value = 3148.6667
$ws.Range("L94").Value = 3148.6667
$ws.Range("N94").Value = -4050.6667
$ws.Range("H132").Value = 25592.176
$ws.Range("I132").Value = 1804.4667
$ws.Range("K132").Value = 5413.4001
$ws.Range("M132").Value = -2883.4001
$ws.Range("H134").Value = 3772.5833
$ws.Range("I134").Value = 4141.5557
$ws.Range("K134").Value = 12424.6671
$ws.Range("M134").Value = -9889.667099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 435651.56
$ws.Range("J5").Value = 1251998.9
$ws.Range("L5").Value = 3755996.7
$ws.Range("N5").Value = -3756220.7
$ws.Range("H92").Value = 362.5
$ws.Range("I92").Value = 295
$ws.Range("J92").Value = 497.5
$ws.Range("K92").Value = 885
$ws.Range("L92").Value = 1492.5
$ws.Range("M92").Value = 363
$ws.Range("N92").Value = -3988.5
$ws.Range("H113").Value = 28249.5
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("H135").Value = 435651.56
$ws.Range("J135").Value = 1251998.9
$ws.Range("L135").Value = 11267990.1
$ws.Range("N135").Value = -11273060.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 14671.667
$ws.Range("I97").Value = 14671.667
$ws.Range("K97").Value = 14671.667
$ws.Range("M97").Value = -14175.667
$ws.Range("H122").Value = 5640.9556
$ws.Range("I122").Value = 3620.8057
$ws.Range("J122").Value = 13721.556
$ws.Range("K122").Value = 10862.4171
$ws.Range("L122").Value = 41164.66800000001
$ws.Range("M122").Value = -8412.417099999999
$ws.Range("N122").Value = -46064.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16419.244
$ws.Range("I7").Value = 24722.422
$ws.Range("K7").Value = 24722.422
$ws.Range("M7").Value = -24610.422
$ws.Range("H16").Value = 5102.0625
$ws.Range("I16").Value = 5422.364
$ws.Range("K16").Value = 5422.364
$ws.Range("M16").Value = -5252.364
$ws.Range("H122").Value = 4599.269
$ws.Range("I122").Value = 4167.75
$ws.Range("J122").Value = 5289.7
$ws.Range("K122").Value = 12503.25
$ws.Range("L122").Value = 15869.1
$ws.Range("M122").Value = -10053.25
$ws.Range("N122").Value = -20769.1
$ws.Range("H126").Value = 16419.244
$ws.Range("I126").Value = 24722.422
$ws.Range("K126").Value = 74167.266
$ws.Range("M126").Value = -71697.266
$ws.Range("H136").Value = 5222.731
$ws.Range("I136").Value = 2623.4119
$ws.Range("K136").Value = 7870.2357
$ws.Range("M136").Value = -5320.2357

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11878
$ws.Range("I81").Value = 17299.5
$ws.Range("J81").Value = 3745.75
$ws.Range("K81").Value = 34599
$ws.Range("L81").Value = 7491.5
$ws.Range("M81").Value = -33538
$ws.Range("N81").Value = -9613.5
$ws.Range("H84").Value = 11878
$ws.Range("I84").Value = 17299.5
$ws.Range("J84").Value = 3745.75
$ws.Range("K84").Value = 172995
$ws.Range("L84").Value = 37457.5
$ws.Range("M84").Value = -167691
$ws.Range("N84").Value = -48065.5
$ws.Range("H132").Value = 17857.459
$ws.Range("I132").Value = 24906.041
$ws.Range("J132").Value = 4844.6924
$ws.Range("K132").Value = 74718.12300000001
$ws.Range("L132").Value = 14534.0772
$ws.Range("M132").Value = -72188.12300000001
$ws.Range("N132").Value = -19594.0772
